$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column after B (becomes column C), shifting old C.. to D..
$ws.Columns("C").Insert()

# Insert a new blank column after the time-of-day column (old G, now H),
# shifting the old I.. columns to K.. (so the ID column ends up at K)
$ws.Columns("I").Insert()

# The insert above copied column B's (green) formatting into the new column C
# for every row that had data in B. Rows 2-4 shouldn't carry that formatting
# at all, so clear them completely.
$ws.Range("C2:C4").Clear()

# Row 5 gets a new header label in the new column C, keeping the format that
# was copied in from column B (date format / green fill).
$ws.Range("C5").Value = "Task/Event number"

# Rows 6-25 become a "Task/Event number" tally column: give them a plain
# integer display format (this produces the new cellXfs entry), then fill in
# the numbers 1-18 for rows 7-24 (row 6 and row 25 stay blank but formatted).
$ws.Range("C6:C25").NumberFormat = "0"
for ($i = 0; $i -lt 18; $i++) {
    $ws.Cells.Item(7 + $i, 3).Value = $i + 1
}

# Tidy up the view state to match where the author left the selection.
$ws.Range("J3").Select()
$excel.ActiveWindow.ScrollColumn = 5
